$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.085.83"
$ws.Range("E2").Value = "  -7.34%  "
$ws.Range("D3").Value = "3.293.53"
$ws.Range("E3").Value = "  -4.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.73"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.78"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.26%  "
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").Value = "3.294.13"
$ws.Range("E8").Value = "  -4.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.465"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.33"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.02%  "
$ws.Range("E11").Value = "  -5.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.368"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.52%  "
$ws.Range("D13").Value = "3.855.39"
$ws.Range("E13").Value = "  -5.00%  "
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "3.295.10"
$ws.Range("E15").Value = "  -4.95%  "
$ws.Range("E16").Value = "  -5.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "23.90"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.37%  "
$ws.Range("D18").Value = "59.216.82"
$ws.Range("E18").Value = "  -7.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.60"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.17"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.85"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -10.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "348.79"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -9.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.550"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.87%  "
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").Value = "3.423.74"
$ws.Range("E25").Value = "  -4.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "68.07"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -8.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000108"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.21"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.44"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.72"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.84%  "
$ws.Range("E32").Value = "  -3.23%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.07"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -6.20%  "
$ws.Range("D35").Value = "3.317.72"
$ws.Range("E35").Value = "  -4.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.55"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.24"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.74"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.47"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "158.20"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0741"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.997"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.34"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.37%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.26"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.07%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.737"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -6.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.16"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.63"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.52"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.68"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("B50").Value = "LidoDAOToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.35"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +13.91%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.46"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +6.29%  "
